# Updates to PID trend stats + chart (current-day reporting window) per commit:
# "updates to graph script to compute the current day. updates to pid trend stats"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated PID trend stats (Sheet1!C2:E8 underlie the chart series) ---
$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 16

$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 7

$ws.Range("C7").Value = 38
$ws.Range("D7").Value = 35
$ws.Range("E7").Value = 37

$ws.Range("C8").Value = 84
$ws.Range("D8").Value = 72
$ws.Range("E8").Value = 79

# --- Chart: give it a title showing the reporting window ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.HasTitle = $true
$chart.ChartTitle.Text = "25-Oct to 7-Nov"

# --- Chart grew taller to accommodate the wider reporting window ---
$co.Height = 519.75

# --- Selection / view moved off the chart and zoom pinned to 100% ---
$ws.Range("D16").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
